$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 109, shifting existing rows 109-125 down to 110-126.
$ws.Range("A109").EntireRow.Insert()

# Populate the newly inserted row 109 with the new weekly record.
$ws.Range("A109").Value = 10
$ws.Range("B109").Value = "Vega Modelo de Temuco"
$ws.Range("C109").Value = "La Araucanía"
$ws.Range("D109").Value = 45124
$ws.Range("E109").Value = 9
$ws.Range("F109").Value = 100112010
$ws.Range("G109").Value = "Achicoria"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 300
$ws.Range("K109").Value = 8000
$ws.Range("L109").Value = 8000
$ws.Range("M109").Value = 8000
$ws.Range("N109").Value = '$/caja 18 unidades'
$ws.Range("O109").Value = "Región Metropolitana"
$ws.Range("P109").Value = 444
$ws.Range("Q109").Value = 18
$ws.Range("R109").Value = "Hortaliza"
